$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (old data rows that are removed in the new version)
$ws.Range("A6:J9").Delete() | Out-Null

# Update the remaining data rows (2-5) with the new values
$data = @(
    @(99888, 85, 72, 74, 58, 74, 82, 70, 24, 33),
    @(99887, 82, 100, 1, 85, 1, 100, 50, 95, 38),
    @(99886, 63, 37, 26, 68, 39, 43, 44, 58, 12),
    @(99885, 32, 73, 77, 76, 28, 52, 22, 49, 63)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

# Update the selection to match the new state: active cell A6, selection A6:XFD13
$ws.Range("A6:XFD13").Select() | Out-Null
